{"js": "// The \"Requisitos\" section ended with:\n//   \"LOB1039: F\u00edsica Experimental III (Requisito fraco)\"\n//   <blank paragraph>\n//   \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n//   \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github\n//    pages. Original theme under Creative Commons Attribution\"\n//   <blank paragraph>\n//   <page-break paragraph>\n// The site-footer boilerplate (the blank line plus the two text\n// paragraphs that follow the requirement line) is removed, leaving the\n// requirement line directly followed by the remaining blank paragraph and\n// the page-break paragraph.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\nconst anchorText = \"LOB1039: F\u00edsica Experimental III (Requisito fraco)\";\nconst removeTexts = [\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n];\n\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === anchorText) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex === -1) {\n  throw new Error(\"Could not locate the 'LOB1039' requirement paragraph\");\n}\n\nconst toDelete = [];\n\n// The blank paragraph right after the requirement line.\nif (items[anchorIndex + 1] && items[anchorIndex + 1].text === \"\") {\n  toDelete.push(anchorIndex + 1);\n}\n\n// The \"Ver no Jupiter...\" / \"\u00a9 2020...\" paragraphs that follow it.\nfor (let i = anchorIndex + 2; i < items.length; i++) {\n  if (removeTexts.indexOf(items[i].text) !== -1) {\n    toDelete.push(i);\n  }\n}\n\n// Delete from the highest index down so earlier indices stay valid.\ntoDelete.sort((a, b) => b - a);\nfor (const idx of toDelete) {\n  items[idx].delete();\n}\n\nawait context.sync();\n", "ps1": "# The \"Requisitos\" section ended with:\n#   \"LOB1039: F\u00edsica Experimental III (Requisito fraco)\"\n#   <blank paragraph>\n#   \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n#   \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github\n#    pages. Original theme under Creative Commons Attribution\"\n#   <blank paragraph>\n#   <page-break paragraph>\n# The site-footer boilerplate (the blank line plus the two text paragraphs\n# that follow the requirement line) is removed, leaving the requirement\n# line directly followed by the remaining blank paragraph and the\n# page-break paragraph.\n$d = $word.ActiveDocument\n\n$anchorText = \"LOB1039: F\u00edsica Experimental III (Requisito fraco)\"\n$removeTexts = @(\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n)\n\n$count = $d.Paragraphs.Count\n$anchorIndex = -1\n\nfor ($i = 1; $i -le $count; $i++) {\n  $text = $d.Paragraphs.Item($i).Range.Text.TrimEnd(\"`r\", \"`a\")\n  if ($text -eq $anchorText) {\n    $anchorIndex = $i\n    break\n  }\n}\n\nif ($anchorIndex -eq -1) {\n  throw \"Could not locate the 'LOB1039' requirement paragraph\"\n}\n\n$toDelete = @()\n\n# The blank paragraph right after the requirement line.\nif ($anchorIndex + 1 -le $count) {\n  $nextText = $d.Paragraphs.Item($anchorIndex + 1).Range.Text.TrimEnd(\"`r\", \"`a\")\n  if ($nextText -eq \"\") {\n    $toDelete += ($anchorIndex + 1)\n  }\n}\n\n# The \"Ver no Jupiter...\" / \"\u00a9 2020...\" paragraphs that follow it.\nfor ($i = $anchorIndex + 2; $i -le $count; $i++) {\n  $text = $d.Paragraphs.Item($i).Range.Text.TrimEnd(\"`r\", \"`a\")\n  if ($removeTexts -contains $text) {\n    $toDelete += $i\n  }\n}\n\n# Delete from the highest index down so earlier indices stay valid.\n$sorted = $toDelete | Sort-Object -Descending\n\nforeach ($idx in $sorted) {\n  $d.Paragraphs.Item($idx).Range.Delete()\n}\n"}
